# Auto-generated edit script applying the Ultima_Profits market-data refresh diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21 (ALC)
$ws.Cells.Item(21, 8).Value = 16079.667
$ws.Cells.Item(21, 9).Value = 1567.7142
$ws.Cells.Item(21, 10).Value = 36396.4
$ws.Cells.Item(21, 11).Value = 1567.7142
$ws.Cells.Item(21, 12).Value = 36396.4
$ws.Cells.Item(21, 13).Value = -1099.7142
$ws.Cells.Item(21, 14).Value = -37332.4

# Row 23 (ALC)
$ws.Cells.Item(23, 8).Value = 16079.667
$ws.Cells.Item(23, 9).Value = 1567.7142
$ws.Cells.Item(23, 10).Value = 36396.4
$ws.Cells.Item(23, 11).Value = 1567.7142
$ws.Cells.Item(23, 12).Value = 36396.4
$ws.Cells.Item(23, 13).Value = -1333.7142
$ws.Cells.Item(23, 14).Value = -36864.4

# Row 34 (ALC)
$ws.Cells.Item(34, 8).Value = 5664.4
$ws.Cells.Item(34, 9).Value = 830.5
$ws.Cells.Item(34, 10).Value = 25000
$ws.Cells.Item(34, 11).Value = 830.5
$ws.Cells.Item(34, 12).Value = 25000
$ws.Cells.Item(34, 13).Value = -627.5
$ws.Cells.Item(34, 14).Value = -25406

# Row 36 (ALC)
$ws.Cells.Item(36, 8).Value = 5664.4
$ws.Cells.Item(36, 9).Value = 830.5
$ws.Cells.Item(36, 10).Value = 25000
$ws.Cells.Item(36, 11).Value = 830.5
$ws.Cells.Item(36, 12).Value = 25000
$ws.Cells.Item(36, 13).Value = -115.5
$ws.Cells.Item(36, 14).Value = -26430

# Row 54 (ALC)
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 13).ClearContents()

# Row 64 (ALC)
$ws.Cells.Item(64, 8).Value = 3666185.8
$ws.Cells.Item(64, 9).Value = 5497428.5
$ws.Cells.Item(64, 10).Value = 3700
$ws.Cells.Item(64, 11).Value = 5497428.5
$ws.Cells.Item(64, 12).Value = 3700
$ws.Cells.Item(64, 13).Value = -5497180.5
$ws.Cells.Item(64, 14).Value = -4196

# Row 67 (ALC)
$ws.Cells.Item(67, 8).Value = 3666185.8
$ws.Cells.Item(67, 9).Value = 5497428.5
$ws.Cells.Item(67, 10).Value = 3700
$ws.Cells.Item(67, 11).Value = 5497428.5
$ws.Cells.Item(67, 12).Value = 3700
$ws.Cells.Item(67, 13).Value = -5496570.5
$ws.Cells.Item(67, 14).Value = -5416

# Row 82 (ALC)
$ws.Cells.Item(82, 8).Value = 3334.8
$ws.Cells.Item(82, 9).Value = 1421.1428
$ws.Cells.Item(82, 10).Value = 7800
$ws.Cells.Item(82, 11).Value = 4263.428400000001
$ws.Cells.Item(82, 12).Value = 23400
$ws.Cells.Item(82, 13).Value = -3857.428400000001
$ws.Cells.Item(82, 14).Value = -24212

# Row 85 (ALC)
$ws.Cells.Item(85, 8).Value = 3334.8
$ws.Cells.Item(85, 9).Value = 1421.1428
$ws.Cells.Item(85, 10).Value = 7800
$ws.Cells.Item(85, 11).Value = 4263.428400000001
$ws.Cells.Item(85, 12).Value = 23400
$ws.Cells.Item(85, 13).Value = -2859.428400000001
$ws.Cells.Item(85, 14).Value = -26208

# Row 138 (ALC)
$ws.Cells.Item(138, 8).Value = 1781.5555
$ws.Cells.Item(138, 9).Value = 1393.6111
$ws.Cells.Item(138, 10).Value = 3333.3333
$ws.Cells.Item(138, 11).Value = 4180.8333
$ws.Cells.Item(138, 12).Value = 9999.999899999999
$ws.Cells.Item(138, 13).Value = 959.1666999999998
$ws.Cells.Item(138, 14).Value = -20279.9999

# Row 139 (ALC)
$ws.Cells.Item(139, 8).Value = 35354.5

$ws = $wb.Worksheets.Item("ARM")
# Row 51 (ARM)
$ws.Cells.Item(51, 8).Value = 20000
$ws.Cells.Item(51, 10).Value = 20000
$ws.Cells.Item(51, 12).Value = 20000
$ws.Cells.Item(51, 14).Value = -21512

# Row 102 (ARM)
$ws.Cells.Item(102, 8).Value = 3250
$ws.Cells.Item(102, 9).Value = 3250
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 3250
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -1628
$ws.Cells.Item(102, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 54 (BSM)
$ws.Cells.Item(54, 8).Value = 1541.5
$ws.Cells.Item(54, 9).Value = 1541.5
$ws.Cells.Item(54, 11).Value = 1541.5
$ws.Cells.Item(54, 13).Value = -1057.5

$ws = $wb.Worksheets.Item("CRP")
# Row 26 (CRP)
$ws.Cells.Item(26, 8).Value = 12234.777
$ws.Cells.Item(26, 9).Value = 1675
$ws.Cells.Item(26, 11).Value = 1675
$ws.Cells.Item(26, 13).Value = -1388

# Row 31 (CRP)
$ws.Cells.Item(31, 8).Value = 6670443.5
$ws.Cells.Item(31, 9).Value = 3880.2368
$ws.Cells.Item(31, 10).Value = 27781226
$ws.Cells.Item(31, 11).Value = 3880.2368
$ws.Cells.Item(31, 12).Value = 27781226
$ws.Cells.Item(31, 13).Value = -3585.2368
$ws.Cells.Item(31, 14).Value = -27781816

# Row 34 (CRP)
$ws.Cells.Item(34, 8).Value = 6670443.5
$ws.Cells.Item(34, 9).Value = 3880.2368
$ws.Cells.Item(34, 10).Value = 27781226
$ws.Cells.Item(34, 11).Value = 3880.2368
$ws.Cells.Item(34, 12).Value = 27781226
$ws.Cells.Item(34, 13).Value = -3678.2368
$ws.Cells.Item(34, 14).Value = -27781630

# Row 44 (CRP)
$ws.Cells.Item(44, 8).Value = 101200
$ws.Cells.Item(44, 9).Value = 170000
$ws.Cells.Item(44, 10).Value = 32400
$ws.Cells.Item(44, 11).Value = 170000
$ws.Cells.Item(44, 12).Value = 32400
$ws.Cells.Item(44, 13).Value = -169558
$ws.Cells.Item(44, 14).Value = -33284

# Row 50 (CRP)
$ws.Cells.Item(50, 8).Value = 9418.5
$ws.Cells.Item(50, 10).Value = 9418.5
$ws.Cells.Item(50, 12).Value = 9418.5
$ws.Cells.Item(50, 14).Value = -10668.5

# Row 56 (CRP)
$ws.Cells.Item(56, 8).Value = 5000
$ws.Cells.Item(56, 9).Value = 5000
$ws.Cells.Item(56, 11).Value = 5000
$ws.Cells.Item(56, 13).Value = -4155

# Row 94 (CRP)
$ws.Cells.Item(94, 8).Value = 3871.818
$ws.Cells.Item(94, 9).Value = 2784
$ws.Cells.Item(94, 10).Value = 4493.4287
$ws.Cells.Item(94, 11).Value = 2784
$ws.Cells.Item(94, 12).Value = 4493.4287
$ws.Cells.Item(94, 13).Value = -2333
$ws.Cells.Item(94, 14).Value = -5395.4287

# Row 105 (CRP)
$ws.Cells.Item(105, 8).Value = 2176.25
$ws.Cells.Item(105, 9).Value = 2315.7144
$ws.Cells.Item(105, 10).Value = 1200
$ws.Cells.Item(105, 11).Value = 2315.7144
$ws.Cells.Item(105, 12).Value = 1200
$ws.Cells.Item(105, 13).Value = -568.7143999999998
$ws.Cells.Item(105, 14).Value = -4694

# Row 132 (CRP)
$ws.Cells.Item(132, 8).Value = 8334772.5
$ws.Cells.Item(132, 9).Value = 10639379
$ws.Cells.Item(132, 10).Value = 2732.7693
$ws.Cells.Item(132, 11).Value = 31918137
$ws.Cells.Item(132, 12).Value = 8198.3079
$ws.Cells.Item(132, 13).Value = -31915607
$ws.Cells.Item(132, 14).Value = -13258.3079

# Row 140 (CRP)
$ws.Cells.Item(140, 8).Value = 48622.5
$ws.Cells.Item(140, 10).Value = 48622.5
$ws.Cells.Item(140, 12).Value = 48622.5
$ws.Cells.Item(140, 14).Value = -58982.5

$ws = $wb.Worksheets.Item("CUL")
# Row 42 (CUL)
$ws.Cells.Item(42, 8).Value = 3399.2222
$ws.Cells.Item(42, 9).Value = 750
$ws.Cells.Item(42, 10).Value = 4156.143
$ws.Cells.Item(42, 11).Value = 2250
$ws.Cells.Item(42, 12).Value = 12468.429
$ws.Cells.Item(42, 13).Value = -1716
$ws.Cells.Item(42, 14).Value = -13536.429

# Row 55 (CUL)
$ws.Cells.Item(55, 8).Value = 588.46155
$ws.Cells.Item(55, 9).Value = 112.5
$ws.Cells.Item(55, 10).Value = 800
$ws.Cells.Item(55, 11).Value = 337.5
$ws.Cells.Item(55, 12).Value = 2400
$ws.Cells.Item(55, 13).Value = -160.5
$ws.Cells.Item(55, 14).Value = -2754

# Row 122 (CUL)
$ws.Cells.Item(122, 8).Value = 687.4194
$ws.Cells.Item(122, 10).Value = 790.6667
$ws.Cells.Item(122, 12).Value = 7116.0003
$ws.Cells.Item(122, 14).Value = -12016.0003

# Row 131 (CUL)
$ws.Cells.Item(131, 8).Value = 825.46
$ws.Cells.Item(131, 9).Value = 288.33334
$ws.Cells.Item(131, 10).Value = 859.7447
$ws.Cells.Item(131, 11).Value = 865.0000200000001
$ws.Cells.Item(131, 12).Value = 2579.2341
$ws.Cells.Item(131, 13).Value = 4174.99998
$ws.Cells.Item(131, 14).Value = -12659.2341

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (GSM)
$ws.Cells.Item(5, 8).Value = 1300
$ws.Cells.Item(5, 9).Value = 166.66667
$ws.Cells.Item(5, 10).Value = 3000
$ws.Cells.Item(5, 11).Value = 166.66667
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 13).Value = -54.66667000000001
$ws.Cells.Item(5, 14).Value = -3224

# Row 46 (GSM)
$ws.Cells.Item(46, 8).Value = 8863.571
$ws.Cells.Item(46, 9).Value = 5000
$ws.Cells.Item(46, 10).Value = 10409
$ws.Cells.Item(46, 11).Value = 5000
$ws.Cells.Item(46, 12).Value = 10409
$ws.Cells.Item(46, 13).Value = -4844
$ws.Cells.Item(46, 14).Value = -10721

# Row 70 (GSM)
$ws.Cells.Item(70, 8).Value = 37049.832
$ws.Cells.Item(70, 9).Value = 102199.5
$ws.Cells.Item(70, 10).Value = 4475
$ws.Cells.Item(70, 11).Value = 102199.5
$ws.Cells.Item(70, 12).Value = 4475
$ws.Cells.Item(70, 13).Value = -101929.5
$ws.Cells.Item(70, 14).Value = -5015

# Row 73 (GSM)
$ws.Cells.Item(73, 8).Value = 37049.832
$ws.Cells.Item(73, 9).Value = 102199.5
$ws.Cells.Item(73, 10).Value = 4475
$ws.Cells.Item(73, 11).Value = 102199.5
$ws.Cells.Item(73, 12).Value = 4475
$ws.Cells.Item(73, 13).Value = -101263.5
$ws.Cells.Item(73, 14).Value = -6347

$ws = $wb.Worksheets.Item("LTW")
# Row 122 (LTW)
$ws.Cells.Item(122, 8).Value = 4431.6113
$ws.Cells.Item(122, 9).Value = 4590.346
$ws.Cells.Item(122, 10).Value = 4018.9
$ws.Cells.Item(122, 11).Value = 13771.038
$ws.Cells.Item(122, 12).Value = 12056.7
$ws.Cells.Item(122, 13).Value = -11321.038
$ws.Cells.Item(122, 14).Value = -16956.7

$ws = $wb.Worksheets.Item("WVR")
# Row 32 (WVR)
$ws.Cells.Item(32, 8).Value = 2000
$ws.Cells.Item(32, 9).Value = 2000
$ws.Cells.Item(32, 11).Value = 2000
$ws.Cells.Item(32, 13).Value = -1683

# Row 34 (WVR)
$ws.Cells.Item(34, 8).Value = 1926
$ws.Cells.Item(34, 9).Value = 1926
$ws.Cells.Item(34, 11).Value = 1926
$ws.Cells.Item(34, 13).Value = -1723

# Row 39 (WVR)
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).ClearContents()
$ws.Cells.Item(39, 14).ClearContents()

# Row 42 (WVR)
$ws.Cells.Item(42, 8).Value = 27400
$ws.Cells.Item(42, 10).Value = 39800
$ws.Cells.Item(42, 12).Value = 39800
$ws.Cells.Item(42, 14).Value = -40556

# Row 138 (WVR)
$ws.Cells.Item(138, 8).Value = 46286
$ws.Cells.Item(138, 10).Value = 46286
$ws.Cells.Item(138, 12).Value = 46286
$ws.Cells.Item(138, 14).Value = -56566
